$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Existing data occupies rows 2..38 (dates descending in col A, deaths in col
# B). A new, more recent day's figures (date serial 43936 => 2020-04-15,
# 130 deaths) are inserted at the top of the table (row 2), which pushes
# every existing data row down by one (the table grows from A1:B38 to
# A1:B39).
#
# Rather than using Rows.Insert() (which drags the header's bold/centered
# style onto the newly shifted-into row), shift the existing values down
# explicitly, working from the bottom of the table upward, then write the
# new row's values. This preserves each cell's existing per-row formatting,
# and the brand-new row 39 naturally inherits the column's default style
# (col A is styled as a date column; col B has no explicit style), matching
# how the rest of the table already looks.

$lastRow = 38

for ($r = $lastRow; $r -ge 2; $r--) {
    $dateVal = $ws.Cells.Item($r, 1).Value2
    $deathVal = $ws.Cells.Item($r, 2).Value2
    $ws.Cells.Item($r + 1, 1).Value = $dateVal
    $ws.Cells.Item($r + 1, 2).Value = $deathVal
}

$ws.Cells.Item(2, 1).Value = 43936
$ws.Cells.Item(2, 2).Value = 130
